$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A98").Value = 7
$ws.Range("B98").Value = 15
$ws.Range("C98").Value = 1.5
$ws.Range("D98").Value = 50
$ws.Range("E98").Value = 82.20999999999999
$ws.Range("F98").Value = 10201
